$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 168, shifting the existing data rows (168:205) down to (169:206).
$ws.Rows.Item(168).Insert()

# Populate the newly inserted row 168 with the new weekly record.
$ws.Cells.Item(168, 1).Value = 5
$ws.Cells.Item(168, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(168, 3).Value = "Maule"
$ws.Cells.Item(168, 4).Value = 44543
$ws.Cells.Item(168, 5).Value = 7
$ws.Cells.Item(168, 6).Value = 100112006
$ws.Cells.Item(168, 7).Value = "Repollo"
$ws.Cells.Item(168, 8).Value = "Crespo record"
$ws.Cells.Item(168, 9).Value = "Primera"
$ws.Cells.Item(168, 10).Value = 5000
$ws.Cells.Item(168, 11).Value = 500
$ws.Cells.Item(168, 12).Value = 500
$ws.Cells.Item(168, 13).Value = 500
$ws.Cells.Item(168, 14).Value = "$/unidad"
$ws.Cells.Item(168, 15).Value = "Región del Maule"
$ws.Cells.Item(168, 16).Value = 500
$ws.Cells.Item(168, 17).Value = 1
$ws.Cells.Item(168, 18).Value = "Hortaliza"
